# edit.ps1 -- reproduce the target diff via Word COM-interop calls.
#
# Summary of content changes:
#  1. " = 113.14 " / "kN" / " T, F"  -> merge into one run, drop proofErr marks
#  2. " = 80 k C, F"                  -> " = 80 kN C, F" (insert N) and move the
#                                        "_GoBack" bookmark to sit right after it
#  3. " = 18,967 lbs " / "C  Design" / " 2: F" -> merge into one run
#  4. " = 21 " / "kN" / " T, F"       -> merge into one run, drop proofErr marks
#  5. "A 20 N force ... a" / "can" / " crushing mechanism..." -> merge into one run
#  6. "Solution: " / "F" / (subscript "can") / " = 148.9 N" -> merge leading
#     "Solution: " + "F" into one run, keep "can" subscript run, drop proofErr marks

$d = $word.ActiveDocument

function Replace-Text($needle, $replacement) {
    $rng = $d.Content
    return $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
}

# 1) "Solution: F[AB] = 113.14 [kN] T, F[AC]" -- collapse the gramStart/gramEnd
#    wrapped "kN" run together with its neighbours.
Replace-Text " = 113.14 kN T, F" " = 113.14 kN T, F" | Out-Null

# 2) "Solution: F[AC] = 80 k C, F[BC]" -- the author typed an "N" after "k" last,
#    so this is also where the _GoBack bookmark (last edit position) now lives.
Replace-Text " = 80 k C, F" " = 80 kN C, F" | Out-Null

$goBackRng = $d.Content
$goBackRng.Find.Execute(" = 80 kN C, F", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $goBackRng.Start + 8   # right after "... = 80 kN"
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 3) "[CD] = 18,967 lbs [C  Design] 2: F[AB]" -- collapse the gramStart/gramEnd
#    wrapped "C  Design" run together with its neighbours.
Replace-Text " = 18,967 lbs C  Design 2: F" " = 18,967 lbs C  Design 2: F" | Out-Null

# 4) "Solution: F[CE] = 21 [kN] T, F[CF]" -- collapse the gramStart/gramEnd
#    wrapped "kN" run together with its neighbours.
Replace-Text " = 21 kN T, F" " = 21 kN T, F" | Out-Null

# 5) "A 20 N force is applied to a [can] crushing mechanism..." -- collapse the
#    spellStart/spellEnd wrapped "can" run together with its neighbours.
Replace-Text "A 20 N force is applied to a can crushing mechanism as shown below. If the distance between points C and D is .1 meters, what are the forces being applied to the can at points B and D?" "A 20 N force is applied to a can crushing mechanism as shown below. If the distance between points C and D is .1 meters, what are the forces being applied to the can at points B and D?" | Out-Null

# 6) "Solution: [F][can] = 148.9 N" -- merge "Solution: " and "F" into one run
#    (this also removes the leading spellStart marker), then restore the
#    subscript formatting on "can" and drop the trailing spellEnd marker.
Replace-Text "Solution: Fcan = 148.9 N" "Solution: Fcan = 148.9 N" | Out-Null

$canRng = $d.Content
$canRng.Find.Execute("Solution: Fcan = 148.9 N", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$canStart = $canRng.Start + 11   # length of "Solution: F"
$canOnly = $d.Range($canStart, $canStart + 3)
$canOnly.Font.Subscript = $true
